# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ) and the
# derived LevePriceNQ/HQ + LeveProfitNQ/HQ columns (H, I, K, and M/N where
# applicable) for leves whose underlying item prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 6853.4287
$ws.Range("I62").Value = 2666.3333
$ws.Range("K62").Value = 2666.3333
$ws.Range("M62").Value = -2042.3333

# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 6853.4287
$ws.Range("I65").Value = 2666.3333
$ws.Range("K65").Value = 13331.6665
$ws.Range("M65").Value = -10211.6665

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 2606501.5
$ws.Range("I76").Value = 4687404
$ws.Range("K76").Value = 4687404
$ws.Range("M76").Value = -4687089

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 2606501.5
$ws.Range("I79").Value = 4687404
$ws.Range("K79").Value = 4687404
$ws.Range("M79").Value = -4686312

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1643.7646
$ws.Range("I98").Value = 1684
$ws.Range("K98").Value = 1684
$ws.Range("M98").Value = -186

# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 21720.092
$ws.Range("I113").Value = 25773.889
$ws.Range("J113").Value = 3478
$ws.Range("K113").Value = 25773.889
$ws.Range("L113").Value = 3478
$ws.Range("M113").Value = -22519.889
$ws.Range("N113").Value = -9986

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1643.7646
$ws.Range("I122").Value = 1684
$ws.Range("K122").Value = 5052
$ws.Range("M122").Value = -2602

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4797.4
$ws.Range("I141").Value = 3965
$ws.Range("J141").Value = 5154.143
$ws.Range("K141").Value = 11895
$ws.Range("L141").Value = 15462.429
$ws.Range("M141").Value = -6715
$ws.Range("N141").Value = -25822.429

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 3087.6338
$ws.Range("I32").Value = 1643.8644
$ws.Range("J32").Value = 10186.167
$ws.Range("K32").Value = 1643.8644
$ws.Range("L32").Value = 10186.167
$ws.Range("M32").Value = -1356.8644
$ws.Range("N32").Value = -10760.167

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Range("H63").Value = 7300
$ws.Range("I63").Value = 7300
$ws.Range("K63").Value = 7300
$ws.Range("M63").Value = -6614

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Range("H66").Value = 7300
$ws.Range("I66").Value = 7300
$ws.Range("K66").Value = 36500
$ws.Range("M66").Value = -33068

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 3159.7856
$ws.Range("J88").Value = 4059.7144
$ws.Range("L88").Value = 4059.7144
$ws.Range("N88").Value = -4871.7144

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 3159.7856
$ws.Range("J91").Value = 4059.7144
$ws.Range("L91").Value = 4059.7144
$ws.Range("N91").Value = -6867.7144

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 654.1818
$ws.Range("J97").Value = 999.5
$ws.Range("L97").Value = 999.5
$ws.Range("N97").Value = -1991.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 687.05884
$ws.Range("I94").Value = 656.2143
$ws.Range("K94").Value = 656.2143
$ws.Range("M94").Value = -205.2143

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 1409
$ws.Range("I99").Value = 994
$ws.Range("K99").Value = 994
$ws.Range("M99").Value = 504

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2849.4517
$ws.Range("I31").Value = 2050.35
$ws.Range("J31").Value = 4302.364
$ws.Range("K31").Value = 2050.35
$ws.Range("L31").Value = 4302.364
$ws.Range("M31").Value = -1755.35
$ws.Range("N31").Value = -4892.364

# Row 32: Daddy's Little Girl | Viper-crested Round Shield
$ws.Range("H32").Value = 3470
$ws.Range("I32").Value = 1293.3334
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 1293.3334
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -977.3334
$ws.Range("N32").Value = -10632

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2849.4517
$ws.Range("I34").Value = 2050.35
$ws.Range("J34").Value = 4302.364
$ws.Range("K34").Value = 2050.35
$ws.Range("L34").Value = 4302.364
$ws.Range("M34").Value = -1848.35
$ws.Range("N34").Value = -4706.364

# Row 51: Greenstone for Greenhorns | Jade Crook
$ws.Range("H51").Value = 30520
$ws.Range("J51").Value = 30520
$ws.Range("L51").Value = 30520
$ws.Range("N51").Value = -31992

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1450670.4
$ws.Range("I58").Value = 2416509
$ws.Range("K58").Value = 2416509
$ws.Range("M58").Value = -2416306

# Row 61: Incant Now, Think Later | Jade Crook
$ws.Range("H61").Value = 30520
$ws.Range("J61").Value = 30520
$ws.Range("L61").Value = 30520
$ws.Range("N61").Value = -31216

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1450670.4
$ws.Range("I136").Value = 2416509
$ws.Range("K136").Value = 7249527
$ws.Range("M136").Value = -7246977

$ws = $wb.Worksheets.Item("CUL")
# Row 115: Mixology | Blood Tomato Juice
$ws.Range("H115").Value = 5332.5557
$ws.Range("J115").Value = 6284.7144
$ws.Range("L115").Value = 18854.1432
$ws.Range("N115").Value = -21204.1432

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 772.48
$ws.Range("J131").Value = 787.5106
$ws.Range("L131").Value = 2362.5318
$ws.Range("N131").Value = -12442.5318

$ws = $wb.Worksheets.Item("GSM")
# Row 47: Wear Your Patriotic Pin | Peridot Choker
$ws.Range("H47").Value = 36015.5
$ws.Range("J47").Value = 36015.5
$ws.Range("L47").Value = 36015.5
$ws.Range("N47").Value = -37151.5

# Row 110: Slimming Down | Stonegold Rapier
$ws.Range("H110").Value = 69998.5
$ws.Range("J110").Value = 69998.5
$ws.Range("L110").Value = 69998.5
$ws.Range("N110").Value = -78178.5

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 1243470.1
$ws.Range("I132").Value = 1673836.1
$ws.Range("J132").Value = 6167.75
$ws.Range("K132").Value = 5021508.300000001
$ws.Range("L132").Value = 18503.25
$ws.Range("M132").Value = -5018978.300000001
$ws.Range("N132").Value = -23563.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 6033.1665
$ws.Range("I16").Value = 6033.1665
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 6033.1665
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -5863.1665
$ws.Range("N16").ClearContents()

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 2804.875
$ws.Range("I22").Value = 5350
$ws.Range("J22").Value = 1956.5
$ws.Range("K22").Value = 5350
$ws.Range("L22").Value = 1956.5
$ws.Range("M22").Value = -5055
$ws.Range("N22").Value = -2546.5

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 2804.875
$ws.Range("I27").Value = 5350
$ws.Range("J27").Value = 1956.5
$ws.Range("K27").Value = 5350
$ws.Range("L27").Value = 1956.5
$ws.Range("M27").Value = -5243
$ws.Range("N27").Value = -2170.5

# Row 32: Men Who Scare Up Goats | Goatskin Targe
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683

# Row 76: Dragoon Drop Rate | Dhalmelskin Breeches of Maiming
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676

# Row 79: Exploiting the Adroit (L) | Dhalmelskin Breeches of Maiming
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 535.6
$ws.Range("I93").Value = 422.42856
$ws.Range("J93").Value = 799.6667
$ws.Range("K93").Value = 422.42856
$ws.Range("L93").Value = 799.6667
$ws.Range("M93").Value = 825.5714399999999
$ws.Range("N93").Value = -3295.6667

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 1349.75
$ws.Range("I100").Value = 1349.75
$ws.Range("K100").Value = 1349.75
$ws.Range("M100").Value = -808.75

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 1713.3334
$ws.Range("I132").Value = 1487
$ws.Range("K132").Value = 4461
$ws.Range("M132").Value = -1931

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 3389.889
$ws.Range("J136").Value = 5956.2856
$ws.Range("L136").Value = 17868.8568
$ws.Range("N136").Value = -22968.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 80: Healing with Flair | Hallowed Ramie Gaskins of Healing
$ws.Range("H80").Value = 79950
$ws.Range("J80").Value = 79950
$ws.Range("L80").Value = 79950
$ws.Range("N80").Value = -81946

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1899
$ws.Range("J81").Value = 1873.75
$ws.Range("L81").Value = 3747.5
$ws.Range("N81").Value = -5869.5

# Row 83: Pants Fit for Battle (L) | Hallowed Ramie Gaskins of Healing
$ws.Range("H83").Value = 79950
$ws.Range("J83").Value = 79950
$ws.Range("L83").Value = 239850
$ws.Range("N83").Value = -249834

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1899
$ws.Range("J84").Value = 1873.75
$ws.Range("L84").Value = 18737.5
$ws.Range("N84").Value = -29345.5

# Row 112: Hair Do No Harm | Iridescent Hat of Healing
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1657.3043
$ws.Range("I132").Value = 1505.409
$ws.Range("K132").Value = 4516.227000000001
$ws.Range("M132").Value = -1986.227000000001
